# Sample Processing - first 4 files chosen and processed for Collegiate Word Ratio (CWR).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 file name first (reuses/replaces the old "Full Text Location" slot) ---
$ws.Range("A2").Value = "Architercture-Messina.pdf"

# --- Header row ---
$ws.Range("A1").Value = "Full Text File Name"
$ws.Range("B1").Value = "Sample File Name"
$ws.Range("C1").Value = "Title"
$ws.Range("D1").Value = "Author"
$ws.Range("E1").Value = "Words Sampled"
$ws.Range("F1").Value = "# Of Collegiate Words"
$ws.Range("G1").Value = "Collegiate Word Ratio"

# --- Row 2 remainder: Architecture ---
$ws.Range("B2").Value = "SampleText.txt"
$ws.Range("C2").Value = "Architecture"
$ws.Range("D2").Value = "Messina, Culler, Pfeiffer…"
$ws.Range("E2").Value = 401
$ws.Range("F2").Value = 19
$ws.Range("G2").Value = 19/401

# --- Row 3: The AGI Containment Problem ---
$ws.Range("A3").Value = "The AGI Containment Problem-Yampolskiy.pdf"
$ws.Range("B3").Value = "SampleText2.txt"
$ws.Range("C3").Value = "The AGI Containment Problem"
$ws.Range("D3").Value = "Babcock, Kramar, Yampolskiy"
$ws.Range("E3").Value = 839
$ws.Range("F3").Value = 39
$ws.Range("G3").Value = 39/839

# --- Row 4: To The "Bestfriend" ... ---
$ws.Range("A4").Value = "To The BestFriend-Collier.txt"
$ws.Range("B4").Value = "SampleText3.txt"
$ws.Range("C4").Value = "To The ""Bestfriend"" I decided I couldn't be friends with anymore"
$ws.Range("D4").Value = "Alexandra Collier"
$ws.Range("E4").Value = 841
$ws.Range("F4").Value = 8
$ws.Range("G4").Value = 8/841

# --- Row 5: Generic anatomy of Escherichia coli 0157h7 outbreaks (sample file named first) ---
$ws.Range("B5").Value = "SampleText4.txt"
$ws.Range("A5").Value = "Generic anatomy of Escherichia coli 0157h7 outbreaks-Eppinger.pdf"
$ws.Range("C5").Value = "Generic anatomy of Escherichia coli 0157h7 outbreaks"
$ws.Range("D5").Value = "Eppinger,Mammel,Leclerc,Ravel,Cebula"
$ws.Range("E5").Value = 649
$ws.Range("F5").Value = 28
$ws.Range("G5").Value = 28/649

# --- Column widths (characters) ---
$ws.Columns.Item(1).ColumnWidth = 42.5
$ws.Columns.Item(3).ColumnWidth = 15.166666666666666
$ws.Columns.Item(4).ColumnWidth = 28.666666666666668
$ws.Columns.Item(5).ColumnWidth = 16.5

# --- View / selection state ---
$ws.Range("D6").Select() | Out-Null
